# Auto-generated edit script applying Typhon_Profits market-data refresh
# Updates currentAveragePrice* / Leve*Price* / LeveProfit* columns (H-N) for the
# rows whose underlying market-board snapshot changed, per the authoritative diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet "ALC" ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H92").Value = 47619640
$ws.Range("I92").Value = 66667224
$ws.Range("K92").Value = 66667224
$ws.Range("M92").Value = -66665976
$ws.Range("H100").Value = 200004200
$ws.Range("I100").Value = 333336640
$ws.Range("J100").Value = 5503
$ws.Range("K100").Value = 333336640
$ws.Range("L100").Value = 5503
$ws.Range("M100").Value = -333336099
$ws.Range("N100").Value = -6585
$ws.Range("H129").Value = 173474.17
$ws.Range("J129").Value = 176505.3
$ws.Range("L129").Value = 529515.8999999999
$ws.Range("N129").Value = -539515.8999999999
$ws.Range("H137").Value = 2249
$ws.Range("I137").Value = 1947.25
$ws.Range("K137").Value = 5841.75
$ws.Range("M137").Value = -3291.75

# ---- Sheet "ARM" ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6897.371
$ws.Range("I32").Value = 4530.6226
$ws.Range("K32").Value = 4530.6226
$ws.Range("M32").Value = -4243.6226
$ws.Range("H45").Value = 3054.9285
$ws.Range("I45").Value = 3299.75
$ws.Range("J45").Value = 2957
$ws.Range("K45").Value = 3299.75
$ws.Range("L45").Value = 2957
$ws.Range("M45").Value = -2922.75
$ws.Range("N45").Value = -3711
$ws.Range("H74").Value = 41668292
$ws.Range("I74").Value = 83333944
$ws.Range("J74").Value = 2639.3333
$ws.Range("K74").Value = 83333944
$ws.Range("L74").Value = 2639.3333
$ws.Range("M74").Value = -83333070
$ws.Range("N74").Value = -4387.3333
$ws.Range("H77").Value = 41668292
$ws.Range("I77").Value = 83333944
$ws.Range("J77").Value = 2639.3333
$ws.Range("K77").Value = 416669720
$ws.Range("L77").Value = 13196.6665
$ws.Range("M77").Value = -416665352
$ws.Range("N77").Value = -21932.6665

# ---- Sheet "BSM" ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1119.6
$ws.Range("I94").Value = 583.2308
$ws.Range("K94").Value = 583.2308
$ws.Range("M94").Value = -132.2308
$ws.Range("H105").Value = 3574112.8
$ws.Range("I105").Value = 2608
$ws.Range("J105").Value = 12502875
$ws.Range("K105").Value = 2608
$ws.Range("L105").Value = 12502875
$ws.Range("M105").Value = -861
$ws.Range("N105").Value = -12506369
$ws.Range("H107").Value = 1183.3572
$ws.Range("I107").Value = 741.55554
$ws.Range("J107").Value = 1978.6
$ws.Range("K107").Value = 741.55554
$ws.Range("L107").Value = 1978.6
$ws.Range("M107").Value = 1178.44446
$ws.Range("N107").Value = -5818.6
$ws.Range("H134").Value = 3145.8975
$ws.Range("I134").Value = 3186.2163
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 9558.6489
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -7023.6489
$ws.Range("N134").Value = -12270

# ---- Sheet "CRP" ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1215.3077
$ws.Range("I16").Value = 916.5
$ws.Range("J16").Value = 1471.4286
$ws.Range("K16").Value = 916.5
$ws.Range("L16").Value = 1471.4286
$ws.Range("M16").Value = -629.5
$ws.Range("N16").Value = -2045.4286
$ws.Range("H22").Value = 394.72726
$ws.Range("I22").Value = 210.66667
$ws.Range("J22").Value = 463.75
$ws.Range("K22").Value = 210.66667
$ws.Range("L22").Value = 463.75
$ws.Range("M22").Value = 139.33333
$ws.Range("N22").Value = -1163.75
$ws.Range("H31").Value = 3446
$ws.Range("I31").Value = 3424.9092
$ws.Range("K31").Value = 3424.9092
$ws.Range("M31").Value = -3129.9092
$ws.Range("H34").Value = 3446
$ws.Range("I34").Value = 3424.9092
$ws.Range("K34").Value = 3424.9092
$ws.Range("M34").Value = -3222.9092
$ws.Range("H105").Value = 1021.8889
$ws.Range("I105").Value = 797.2
$ws.Range("K105").Value = 797.2
$ws.Range("M105").Value = 949.8
$ws.Range("H107").Value = 1801.0526
$ws.Range("I107").Value = 751.8
$ws.Range("K107").Value = 751.8
$ws.Range("M107").Value = 1168.2
$ws.Range("H113").Value = 1215.3077
$ws.Range("I113").Value = 916.5
$ws.Range("J113").Value = 1471.4286
$ws.Range("K113").Value = 916.5
$ws.Range("L113").Value = 1471.4286
$ws.Range("M113").Value = 1253.5
$ws.Range("N113").Value = -5811.4286
$ws.Range("H132").Value = 2262.8108
$ws.Range("I132").Value = 1538.6451
$ws.Range("K132").Value = 4615.9353
$ws.Range("M132").Value = -2085.9353

# ---- Sheet "CUL" ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 899
$ws.Range("J75").Value = 848.5
$ws.Range("L75").Value = 2545.5
$ws.Range("N75").Value = -4541.5
$ws.Range("H78").Value = 899
$ws.Range("J78").Value = 848.5
$ws.Range("L78").Value = 7636.5
$ws.Range("N78").Value = -17620.5
$ws.Range("H105").Value = 8676.333000000001
$ws.Range("J105").Value = 8676.333000000001
$ws.Range("L105").Value = 26028.999
$ws.Range("N105").Value = -31270.999
$ws.Range("H129").Value = 223114.56
$ws.Range("I129").Value = 728.8889
$ws.Range("J129").Value = 366076.78
$ws.Range("K129").Value = 2186.6667
$ws.Range("L129").Value = 1098230.34
$ws.Range("M129").Value = 2813.3333
$ws.Range("N129").Value = -1108230.34
$ws.Range("H131").Value = 722.5
$ws.Range("J131").Value = 725
$ws.Range("L131").Value = 2175
$ws.Range("N131").Value = -12255

# ---- Sheet "GSM" ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 32071.646
$ws.Range("I132").Value = 2199.3635
$ws.Range("K132").Value = 6598.0905
$ws.Range("M132").Value = -4068.0905

# ---- Sheet "LTW" ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1027.9395
$ws.Range("I46").Value = 1000.70966
$ws.Range("K46").Value = 1000.70966
$ws.Range("M46").Value = -812.70966

# ---- Sheet "WVR" ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1528.8
$ws.Range("I100").Value = 1528.8
$ws.Range("K100").Value = 3057.6
$ws.Range("M100").Value = -2516.6
$ws.Range("H122").Value = 892.8205
$ws.Range("I122").Value = 838.04
$ws.Range("J122").Value = 990.6429000000001
$ws.Range("K122").Value = 2514.12
$ws.Range("L122").Value = 2971.9287
$ws.Range("M122").Value = -64.11999999999989
$ws.Range("N122").Value = -7871.9287
$ws.Range("H126").Value = 1086.8
$ws.Range("I126").Value = 1133.5
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 3400.5
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -930.5
$ws.Range("N126").Value = -7640
$ws.Range("H136").Value = 23811974
$ws.Range("I136").Value = 31251142
$ws.Range("K136").Value = 93753426
$ws.Range("M136").Value = -93750876

